# Update crypto price/volume snapshot on Sheet1 (GitHub Actions scrape refresh).
# Values are written with a leading apostrophe so Excel keeps them as literal
# text (matching the original workbook, where D/E columns are text, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'278.15"
$ws.Range("E2").Value = "'0.80%"
$ws.Range("D3").Value = "'27.20"
$ws.Range("E3").Value = "'1.61%"
$ws.Range("D4").Value = "'4.865"
$ws.Range("E4").Value = "'-0.18%"
$ws.Range("D5").Value = "'0.06424"
$ws.Range("E5").Value = "'1.49%"
$ws.Range("D6").Value = "'7.014"
$ws.Range("E6").Value = "'1.25%"
$ws.Range("D7").Value = "'1.195"
$ws.Range("E7").Value = "'-6.71%"
$ws.Range("D8").Value = "'0.8863"
$ws.Range("E8").Value = "'1.32%"
$ws.Range("E9").Value = "'-0.78%"
$ws.Range("D10").Value = "'0.05156"
$ws.Range("E10").Value = "'1.29%"
$ws.Range("D11").Value = "'0.07513"
$ws.Range("E11").Value = "'0.51%"
$ws.Range("D12").Value = "'0.02878"
$ws.Range("D13").Value = "'0.08969"
$ws.Range("D14").Value = "'0.001565"
$ws.Range("E14").Value = "'-0.49%"
$ws.Range("D15").Value = "'0.0006377"
$ws.Range("E15").Value = "'0.86%"
$ws.Range("D16").Value = "'0.006091"
$ws.Range("E16").Value = "'1.17%"
$ws.Range("D17").Value = "'3.477"
$ws.Range("E17").Value = "'0.71%"
$ws.Range("D18").Value = "'3.305"
$ws.Range("E18").Value = "'-0.54%"
$ws.Range("D19").Value = "'2.241"
$ws.Range("E19").Value = "'-1.88%"
$ws.Range("E21").Value = "'0.53%"
$ws.Range("D22").Value = "'3.905"
$ws.Range("E22").Value = "'-0.92%"
$ws.Range("D23").Value = "'0.1518"
$ws.Range("E23").Value = "'10.01%"
$ws.Range("E24").Value = "'0.76%"
$ws.Range("D25").Value = "'0.001174"
$ws.Range("E25").Value = "'0.68%"
$ws.Range("D26").Value = "'0.003882"
$ws.Range("E26").Value = "'-7.86%"
$ws.Range("E28").Value = "'-1.72%"
$ws.Range("E29").Value = "'1.73%"
$ws.Range("D40").Value = "'0.04122"
$ws.Range("E40").Value = "'0.65%"
$ws.Range("D41").Value = "'0.006811"
$ws.Range("E41").Value = "'-2.96%"
$ws.Range("E42").Value = "'0.15%"
$ws.Range("D43").Value = "'0.001910"
$ws.Range("E43").Value = "'-12.83%"
$ws.Range("D44").Value = "'0.01170"
$ws.Range("E44").Value = "'3.78%"
$ws.Range("D45").Value = "'0.00005333"
$ws.Range("E45").Value = "'2.20%"
$ws.Range("D46").Value = "'1.683"
$ws.Range("E46").Value = "'13.25%"
